# Generate Report for Handoff
# Marks the fa5d4b70-21a4-4027-bb63-cb76d20052b0.md file as "Ready for handoff"
# across the Overview / zh-cn / de-de sheets, refreshes the related
# handoff timestamps, and records the "stale handback" error detail for
# each locale sheet.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/86dc5cd34bfafb8a8021b0f5e817aa4ef525fd25/e2e/fa5d4b70-21a4-4027-bb63-cb76d20052b0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f79a9c7b611c5c31d4e35e71fa7e71b3f49e8a7/e2e/fa5d4b70-21a4-4027-bb63-cb76d20052b0.md."

# --- Overview sheet: row 3 is fa5d4b70-21a4-4027-bb63-cb76d20052b0.md ---
$ws_overview.Range("E3").Value = "Ready for handoff"
$ws_overview.Range("F3").Value = "Ready for handoff"
$ws_overview.Range("G3").Value = "2016-08-24 00:46:27"

# --- zh-cn sheet: row 3 is fa5d4b70-21a4-4027-bb63-cb76d20052b0.md ---
$ws_zhcn.Range("C3").Value = "Ready for handoff"
$ws_zhcn.Range("H3").Value = "2016-08-24 00:46:22"
$ws_zhcn.Range("P3").Value = $errorDetail
$ws_zhcn.Columns.Item(16).ColumnWidth = 39.14

# --- de-de sheet: row 3 is fa5d4b70-21a4-4027-bb63-cb76d20052b0.md ---
$ws_dede.Range("C3").Value = "Ready for handoff"
$ws_dede.Range("H3").Value = "2016-08-24 00:46:27"
$ws_dede.Range("P3").Value = $errorDetail
$ws_dede.Columns.Item(16).ColumnWidth = 39.14
